$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 1P" ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D4").Value = 10
$ws1.Range("F4").Value = 29
$ws1.Range("G4").Value = 74.36
$ws1.Range("H4").Value = 7.9

$ws1.Range("D5").Value = 6
$ws1.Range("F5").Value = 15
$ws1.Range("G5").Value = 71.43000000000001
$ws1.Range("H5").Value = 9.199999999999999

# --- Sheet "Estadisticos 2P" ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("D2").Value = 4
$ws2.Range("E2").Value = 0
$ws2.Range("F2").Value = 17
$ws2.Range("G2").Value = 80.95
$ws2.Range("H2").Value = 9.1

$ws2.Range("D3").Value = 4
$ws2.Range("E3").Value = 0
$ws2.Range("F3").Value = 17
$ws2.Range("G3").Value = 80.95
$ws2.Range("H3").Value = 9.1

$ws2.Range("D4").Value = 10
$ws2.Range("E4").Value = 0
$ws2.Range("F4").Value = 29
$ws2.Range("G4").Value = 74.36
$ws2.Range("H4").Value = 7.9

$ws2.Range("D5").Value = 6
$ws2.Range("E5").Value = 0
$ws2.Range("F5").Value = 15
$ws2.Range("G5").Value = 71.43000000000001
$ws2.Range("H5").Value = 9.199999999999999

# --- Sheet "Estadisticos Final" ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("H2").Value = 8.6
$ws3.Range("H3").Value = 8.699999999999999

$ws3.Range("D4").Value = 10
$ws3.Range("F4").Value = 29
$ws3.Range("G4").Value = 74.36

$ws3.Range("D5").Value = 6
$ws3.Range("F5").Value = 15
$ws3.Range("G5").Value = 71.43000000000001
$ws3.Range("H5").Value = 9.300000000000001

# --- Sheet "Rescatables" ---
$ws4 = $wb.Worksheets.Item("Rescatables")
$ws4.Rows.Item(2).Delete()
